$wb = $excel.ActiveWorkbook

# "Generate Report for Archive": the handoff-status report is regenerated,
# so every cell that still said "Ready for handoff" now reads
# "In Translation" (zh-cn / de-de detail sheets, plus their rollup on the
# Overview sheet).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $used.Replace("Ready for handoff", "In Translation") | Out-Null
}

# The Status columns (and the Overview sheet's zh-cn/de-de rollup columns,
# which mirror the same text) were sized with AutoFit for the old, longer
# string. Re-fit them now that the text is shorter.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).AutoFit() | Out-Null
$overview.Columns.Item(6).AutoFit() | Out-Null

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).AutoFit() | Out-Null

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).AutoFit() | Out-Null
